$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'293.93"
$ws.Range("E2").Value = "'1.63%"
$ws.Range("D3").Value = "'31.21"
$ws.Range("E3").Value = "'1.32%"
$ws.Range("D4").Value = "'4.943"
$ws.Range("E4").Value = "'1.11%"
$ws.Range("D5").Value = "'0.07334"
$ws.Range("E5").Value = "'2.35%"
$ws.Range("D6").Value = "'2.295"
$ws.Range("E6").Value = "'25.25%"
$ws.Range("D7").Value = "'7.705"
$ws.Range("E7").Value = "'0.67%"
$ws.Range("E8").Value = "'0.35%"
$ws.Range("D9").Value = "'0.9073"
$ws.Range("E9").Value = "'0.91%"
$ws.Range("E10").Value = "'2.16%"
$ws.Range("D11").Value = "'0.08121"
$ws.Range("E11").Value = "'8.39%"
$ws.Range("D12").Value = "'0.08176"
$ws.Range("E12").Value = "'0.66%"
$ws.Range("D13").Value = "'0.03108"
$ws.Range("E13").Value = "'4.17%"
$ws.Range("D14").Value = "'0.1006"
$ws.Range("E14").Value = "'0.61%"
$ws.Range("D15").Value = "'0.001508"
$ws.Range("E15").Value = "'0.40%"
$ws.Range("D16").Value = "'0.005754"
$ws.Range("E16").Value = "'1.03%"
$ws.Range("D17").Value = "'3.483"
$ws.Range("E17").Value = "'0.46%"
$ws.Range("D18").Value = "'2.080"
$ws.Range("E18").Value = "'-1.24%"
$ws.Range("D19").Value = "'0.3329"
$ws.Range("E19").Value = "'1.57%"
$ws.Range("D20").Value = "'0.1287"
$ws.Range("E20").Value = "'-0.97%"
$ws.Range("D21").Value = "'3.971"
$ws.Range("E21").Value = "'-9.13%"
$ws.Range("D22").Value = "'0.2103"
$ws.Range("E22").Value = "'4.89%"
$ws.Range("D23").Value = "'0.04531"
$ws.Range("E23").Value = "'1.41%"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-0.32%"
$ws.Range("D25").Value = "'0.004342"
$ws.Range("E25").Value = "'8.03%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'3.84%"
$ws.Range("D27").Value = "'0.0003398"
$ws.Range("E27").Value = "'-95.49%"
$ws.Range("D39").Value = "'0.01601"
$ws.Range("E39").Value = "'-2.17%"
$ws.Range("D40").Value = "'0.04431"
$ws.Range("E40").Value = "'2.27%"
$ws.Range("E41").Value = "'0.08%"
$ws.Range("D42").Value = "'0.009089"
$ws.Range("D43").Value = "'0.1327"
$ws.Range("E43").Value = "'1.66%"
$ws.Range("D44").Value = "'0.001922"
$ws.Range("E44").Value = "'-4.38%"
$ws.Range("D45").Value = "'0.009198"
$ws.Range("E45").Value = "'-9.47%"
$ws.Range("D46").Value = "'0.00005964"
$ws.Range("E46").Value = "'1.79%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("E48").Value = "'1.87%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.15%"
